$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6: Advanced Leather Rack - update C6, E6, F6
$ws.Range("C6").Value = 1460223
$ws.Range("E6").Value = 19.72761499999999
$ws.Range("F6").Value = 13.51000155455707

# Row 9: becomes "Vintage Black Lion Weapon Box" (new item inserted)
$ws.Range("A9").Value = "Vintage Black Lion Weapon Box"
$ws.Range("B9").Value = 26458
$ws.Range("C9").Value = 400019
$ws.Range("D9").Value = 523332
$ws.Range("E9").Value = 4.481320000000001
$ws.Range("F9").Value = 11.20276786852625
$ws.Range("I9").Value = 103.885

# Row 10: becomes old row 9's data "Ghostly Infusion"
$ws.Range("A10").Value = "Ghostly Infusion"
$ws.Range("B10").Value = 24207
$ws.Range("C10").Value = 624998
$ws.Range("D10").Value = 784994
$ws.Range("E10").Value = 4.224690000000002
$ws.Range("F10").Value = 6.759525630482022
$ws.Range("I10").Value = 29.382
$ws.Range("K10").Value = 1

# Row 11: becomes old row 10's data "Zojja''s Berserker Insignia" but with slightly updated D, E, F, I, L, K values
$ws.Range("A11").Value = "Zojja''s Berserker Insignia"
$ws.Range("B11").Value = 20887
$ws.Range("C11").Value = 179140
$ws.Range("D11").Value = 244480
$ws.Range("E11").Value = 2.8668
$ws.Range("F11").Value = 16.00312604666741
$ws.Range("I11").Value = 27.681
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 1
